$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New "Save" header in H1 - copy formatting from the existing header cell (G1)
$ws.Range("H1").Value = "Save"
$ws.Range("G1").Copy()
$ws.Range("H1").PasteSpecial(-4122)

# New Save values for each data row
$ws.Range("H2").Value = 1
$ws.Range("H3").Value = 1
